# Literature survey - feature selection and stuff
#
# Adds a new row (41) to the schedule sheet describing a literature-survey
# task on feature selection, and moves the current selection/scroll
# position to reflect the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content -------------------------------------------------
# Columns A/B/C reuse existing shared strings (same chapter/task-type as
# the other "Literature review" rows 39/40).
$ws.Range("A41").Value = "נספחים"
$ws.Range("B41").Value = "תאורטי"
$ws.Range("C41").Value = "סקר ספרות"

# New shared strings must be introduced in the same order the original
# workbook recorded them (I -> E -> D) so that they land on the expected
# shared-string table indices.
$ws.Range("I41").Value = "נושא יחסית גנרי (או שלא הבנתי טוב)`nפוסט נחמד בנושא https://www.kaggle.com/residentmario/automated-feature-selection-with-sklearn"
$ws.Range("E41").Value = "7/30/2020 16:00 PM"
$ws.Range("F41").Value = 44042.816666666666
$ws.Range("D41").Value = "לקרוא, לחקור, לסכם, להשוות, לשאול שאלות, למצוא תשובות על בחירת פיצ'רים"
$ws.Range("G41").Value = 0.5
$ws.Range("H41").Value = 0.5

# D/I hold long free-text notes, like the other rows in the table - wrap them.
$ws.Range("D41").WrapText = $true
$ws.Range("I41").WrapText = $true

# Match the row height Excel computed for the wrapped text in the source file.
$ws.Rows.Item(41).RowHeight = 43.5

# --- View state --------------------------------------------------------
# Reflect the newly-added row in the window scroll position / selection.
$null = $ws.Range("D42").Select()
